$d = $word.ActiveDocument

# Insert a new centered paragraph containing "{image-timestamp}" right
# after the image paragraph and before the "{caption-text}" paragraph,
# by inserting a new paragraph mark immediately before the caption
# paragraph's range (keeps the new run free of inherited rPr, e.g. the
# noProof flag carried by the image's run).
$captionPara = $d.Paragraphs(2)
$captionPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs(2)
$newPara.Range.Text = "{image-timestamp}"
$newPara.Format.Alignment = 1
